# Update countries data (Pais sheet) - rows reshuffled per updated COVID stats
# and the "Datos actualizados" timestamp, per commit "Update countries & provincias Spain"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 16:50"
# Row 4
$ws.Range("B4").Value = 165392
$ws.Range("C4").Value = 1604
$ws.Range("D4").Value = 5544
$ws.Range("E4").Value = 156666
$ws.Range("F4").Value = 3535
$ws.Range("G4").Value = 41
$ws.Range("H4").Value = 3182
# Row 6
$ws.Range("E6").Value = 66889
$ws.Range("G6").Value = 553
$ws.Range("H6").Value = 8269
# Row 8
$ws.Range("B8").Value = 68180
$ws.Range("C8").Value = 1295
$ws.Range("E8").Value = 51674
# Row 11
$ws.Range("E11").Value = 20198
$ws.Range("G11").Value = 400
$ws.Range("H11").Value = 1808
# Row 31
$ws.Range("E31").Value = 1947
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 78
# Row 46
$ws.Range("A46").Value = "Republica Dominicana"
$ws.Range("B46").Value = 1109
$ws.Range("C46").Value = 208
$ws.Range("D46").Value = 5
$ws.Range("E46").Value = 1053
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 9
$ws.Range("H46").Value = 51
# Row 47
$ws.Range("A47").Value = "Mexico"
$ws.Range("B47").Value = 1094
$ws.Range("C47").Value = 101
$ws.Range("D47").Value = 35
$ws.Range("E47").Value = 1031
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 28
# Row 48
$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 1075
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 9
$ws.Range("E48").Value = 1039
$ws.Range("F48").Value = 43
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 27
# Row 49
$ws.Range("A49").Value = "Argentina"
$ws.Range("B49").Value = 966
$ws.Range("C49").Value = 146
$ws.Range("D49").Value = 228
$ws.Range("E49").Value = 712
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 26
# Row 50
$ws.Range("A50").Value = "Peru"
$ws.Range("B50").Value = 950
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 53
$ws.Range("E50").Value = 873
$ws.Range("F50").Value = 49
$ws.Range("H50").Value = 24
# Row 51
$ws.Range("A51").Value = "Singapur"
$ws.Range("B51").Value = 926
$ws.Range("C51").Value = 47
$ws.Range("D51").Value = 240
$ws.Range("E51").Value = 683
$ws.Range("F51").Value = 22
$ws.Range("H51").Value = 3
# Row 78
$ws.Range("A78").Value = "Moldavia"
$ws.Range("B78").Value = 353
$ws.Range("C78").Value = 55
$ws.Range("D78").Value = 18
$ws.Range("E78").Value = 333
$ws.Range("F78").Value = 44
$ws.Range("G78").Value = 0
# Row 79
$ws.Range("A79").Value = "Kazajistan"
$ws.Range("B79").Value = 336
$ws.Range("C79").Value = 34
$ws.Range("D79").Value = 22
$ws.Range("E79").Value = 312
$ws.Range("F79").Value = 6
$ws.Range("G79").Value = 1
# Row 80
$ws.Range("A80").Value = "Costa Rica"
$ws.Range("B80").Value = 330
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 4
$ws.Range("E80").Value = 324
$ws.Range("F80").Value = 7
$ws.Range("G80").Value = 0
# Row 81
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 329
$ws.Range("C81").Value = 44
$ws.Range("D81").Value = 12
$ws.Range("E81").Value = 308
$ws.Range("F81").Value = 1
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 9
# Row 82
$ws.Range("A82").Value = "Taiwan"
$ws.Range("B82").Value = 322
$ws.Range("C82").Value = 16
$ws.Range("D82").Value = 39
$ws.Range("E82").Value = 278
$ws.Range("F82").Value = 0
$ws.Range("H82").Value = 5
# Row 83
$ws.Range("A83").Value = "Uruguay"
$ws.Range("B83").Value = 320
$ws.Range("D83").Value = 25
$ws.Range("E83").Value = 294
$ws.Range("F83").Value = 9
$ws.Range("H83").Value = 1
# Row 87
$ws.Range("A87").Value = "Reunion"
$ws.Range("B87").Value = 247
$ws.Range("C87").Value = 23
$ws.Range("D87").Value = 1
$ws.Range("F87").Value = 4
$ws.Range("H87").Value = 0
# Row 88
$ws.Range("A88").Value = "Burkina Faso"
$ws.Range("B88").Value = 246
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 31
$ws.Range("E88").Value = 203
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 12
# Row 89
$ws.Range("A89").Value = "Albania"
$ws.Range("B89").Value = 243
$ws.Range("C89").Value = 20
$ws.Range("D89").Value = 52
$ws.Range("E89").Value = 178
$ws.Range("F89").Value = 8
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 13
# Row 90
$ws.Range("A90").Value = "Republica de Chipre"
$ws.Range("D90").Value = 22
$ws.Range("E90").Value = 201
$ws.Range("F90").Value = 3
$ws.Range("H90").Value = 7
# Row 91
$ws.Range("A91").Value = "San Marino"
$ws.Range("B91").Value = 230
$ws.Range("D91").Value = 13
$ws.Range("E91").Value = 192
$ws.Range("F91").Value = 16
$ws.Range("H91").Value = 25
# Row 104
$ws.Range("E104").Value = 139
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 4
# Row 161
$ws.Range("A161").Value = "Guinea Ecuatorial"
$ws.Range("D161").Value = 1
$ws.Range("H161").Value = 0
# Row 162
$ws.Range("A162").Value = "Islas Caimanes"
$ws.Range("D162").Value = 0
$ws.Range("H162").Value = 1
# Row 169
$ws.Range("A169").Value = "Laos"
$ws.Range("C169").Value = 1
# Row 170
$ws.Range("A170").Value = "Granada"
$ws.Range("C170").Value = 0
# Row 173
$ws.Range("A173").Value = "Mozambique"
# Row 174
$ws.Range("A174").Value = "Libia"
# Row 175
$ws.Range("A175").Value = "Guinea-Bisau"
# Row 177
$ws.Range("A177").Value = "Surinam"
# Row 182
$ws.Range("A182").Value = "Sudan"
$ws.Range("C182").Value = 1
# Row 183
$ws.Range("A183").Value = "Angola"
$ws.Range("C183").Value = 0
# Row 186
$ws.Range("A186").Value = "Benin"
$ws.Range("D186").Value = 1
$ws.Range("H186").Value = 0
# Row 187
$ws.Range("A187").Value = "Cabo Verde"
$ws.Range("D187").Value = 0
$ws.Range("H187").Value = 1
# Row 190
$ws.Range("A190").Value = "Islas Turcas y Caicos"
# Row 191
$ws.Range("A191").Value = "Fiyi"
# Row 192
$ws.Range("A192").Value = "Montserrat"
# Row 195
$ws.Range("A195").Value = "Nicaragua"
# Row 196
$ws.Range("A196").Value = "Gambia"
# Row 197
$ws.Range("A197").Value = "Republica de Africa Central"
# Row 198
$ws.Range("A198").Value = "Botsuana"
# Row 199
$ws.Range("A199").Value = "Belice"
$ws.Range("C199").Value = 0
# Row 200
$ws.Range("A200").Value = "Islas Virgenes Britanicas"
$ws.Range("C200").Value = 1
# Row 201
$ws.Range("A201").Value = "Liberia"
# Row 204
$ws.Range("A204").Value = "Papua Nueva Guinea"
$ws.Range("C204").Value = 0
# Row 205
$ws.Range("A205").Value = "Timor Oriental"
# Row 206
$ws.Range("A206").Value = "Sierra Leona"
$ws.Range("C206").Value = 1
